# NN-368 - Paymaart - Admin Web - Insights of Merchant Registration BDD fix
# Update the sample BDD testing row with a fresh Paymaart ID + phone number,
# tighten up rows 1-3 to the new row height, and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 holds the "BDD TESTING" sample record - refresh the generated
# Paymaart ID and phone number used by the test fixture.
$ws.Range("C2").Value = "CMR15019448"
$ws.Range("D2").Value = "265 46 419 2496"

# Shrink the header/data rows from 15.75pt to 13.8pt.
$ws.Range("A1:F3").RowHeight = 13.8

# Move the active selection from D2 to E3.
$ws.Range("E3").Select() | Out-Null
